$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh customer data: append an "s" to the test names, bump the pin
# codes to new unused values, and tweak the two test-email local parts so
# guru99 sees these as "new" customers (the app can't delete previously
# created ones yet, so re-running with identical data makes tests fail).
$ws.Range("A2").Value = "alexiss"
$ws.Range("G2").Value = '"222222"'
$ws.Range("I2").Value = "artahAlsd+1@gmail.com"

$ws.Range("A3").Value = "Dinas"
$ws.Range("G3").Value = '"333333"'
$ws.Range("I3").Value = "stsADDdDN+Din@gmail.com"

# Move the active selection from I7 to D7.
$ws.Range("D7").Select()

# Narrow the mobile column and widen the email column.
$ws.Columns.Item(8).ColumnWidth = 15.71
$ws.Columns.Item(9).ColumnWidth = 27.75
